$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 277, shifting existing rows 277-305 down to 278-306
$ws.Rows.Item(277).Insert()

# Populate the new row 277 with the new record's data
$ws.Range("A277").Value = 10
$ws.Range("B277").Value = "Vega Modelo de Temuco"
$ws.Range("C277").Value = "La Araucanía"
$ws.Range("D277").Value = 44918
$ws.Range("E277").Value = 9
$ws.Range("F277").Value = 100112052
$ws.Range("G277").Value = "Albahaca"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 125
$ws.Range("K277").Value = 8000
$ws.Range("L277").Value = 8000
$ws.Range("M277").Value = 8000
$ws.Range("N277").Value = "$/paquete"
$ws.Range("O277").Value = "Región Metropolitana"
$ws.Range("P277").Value = 8000
$ws.Range("Q277").Value = 1
$ws.Range("R277").Value = "Hortaliza"
